# Update the "Fixed" date shown in the Date/Footer placeholder from
# 4/24/2022 to 4/30/2022 across the slide master and every slide layout
# (cleanup of supervised learning notebook).

$p = $ppt.ActivePresentation

$oldDate = "4/24/2022"
$newDate = "4/30/2022"

# Collect the slide master plus every custom layout that hangs off it -
# this is where the Date Placeholder shapes actually live.
$targets = @()
$targets += $p.SlideMaster
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $targets += $p.SlideMaster.CustomLayouts.Item($i)
}

foreach ($target in $targets) {
    for ($j = 1; $j -le $target.Shapes.Count; $j++) {
        $shape = $target.Shapes.Item($j)
        if ($shape.Name -like "Date Placeholder*" -and $shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}
